$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "230696 - Carlos Jose Todero Peixoto" row (row 13, B/C only).
# Everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# Update cell contents to their new (shifted) values.
$ws.Range("B10").Value = "230696 - Carlos José Todero Peixoto"
$ws.Range("C10").Value = "230696 - Carlos José Todero Peixoto"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "230696 - Carlos José Todero Peixoto"
$ws.Range("C18").Value = "230696 - Carlos José Todero Peixoto"

$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
